$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '30.807.77'
Set-TextValue $ws.Range("E2") '  -0.79%  '
Set-TextValue $ws.Range("D3") '1.942.39'
Set-TextValue $ws.Range("E3") '  -0.76%  '
Set-TextValue $ws.Range("E4") '  -0.19%  '
Set-TextValue $ws.Range("D5") '241.94'
Set-TextValue $ws.Range("E5") '  -1.86%  '
Set-TextValue $ws.Range("E6") '  -0.15%  '
Set-TextValue $ws.Range("D7") '0.4888'
Set-TextValue $ws.Range("E7") '  -0.35%  '
Set-TextValue $ws.Range("D8") '0.2955'
Set-TextValue $ws.Range("E8") '  -0.35%  '
Set-TextValue $ws.Range("D9") '0.06899'
Set-TextValue $ws.Range("E9") '  +1.23%  '
Set-TextValue $ws.Range("D10") '19.44'
Set-TextValue $ws.Range("E10") '  +2.34%  '
Set-TextValue $ws.Range("D11") '106.41'
Set-TextValue $ws.Range("E11") '  +0.15%  '
Set-TextValue $ws.Range("D12") '1.943.53'
Set-TextValue $ws.Range("E12") '  -0.40%  '
Set-TextValue $ws.Range("D13") '0.07719'
Set-TextValue $ws.Range("E13") '  -0.32%  '
Set-TextValue $ws.Range("D14") '5.344'
Set-TextValue $ws.Range("E14") '  -1.15%  '
Set-TextValue $ws.Range("D15") '0.6993'
Set-TextValue $ws.Range("E15") '  -1.38%  '
Set-TextValue $ws.Range("D16") '276.02'
Set-TextValue $ws.Range("E16") '  -2.53%  '
Set-TextValue $ws.Range("D17") '30.811.89'
Set-TextValue $ws.Range("E17") '  -0.90%  '
Set-TextValue $ws.Range("D18") '0.000007725'
Set-TextValue $ws.Range("E18") '  -0.55%  '
Set-TextValue $ws.Range("D19") '13.12'
Set-TextValue $ws.Range("E19") '  -0.82%  '
Set-TextValue $ws.Range("B20") 'Dai'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D20") '0.9999'
Set-TextValue $ws.Range("E20") '  -0.15%  '
Set-TextValue $ws.Range("B21") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D21") '2.192.39'
Set-TextValue $ws.Range("E21") '  -0.10%  '
Set-TextValue $ws.Range("D22") '5.456'
Set-TextValue $ws.Range("E22") '  -1.68%  '
Set-TextValue $ws.Range("D23") '0.9995'
Set-TextValue $ws.Range("E23") '  -0.23%  '
Set-TextValue $ws.Range("E24") '  -1.20%  '
Set-TextValue $ws.Range("D25") '9.717'
Set-TextValue $ws.Range("E25") '  -2.48%  '
Set-TextValue $ws.Range("D26") '167.67'
Set-TextValue $ws.Range("D27") '19.65'
Set-TextValue $ws.Range("E27") '  -1.60%  '
Set-TextValue $ws.Range("D28") '2.154'
Set-TextValue $ws.Range("E28") '  -1.62%  '
Set-TextValue $ws.Range("E29") '  -1.14%  '
Set-TextValue $ws.Range("D30") '1.392'
Set-TextValue $ws.Range("E30") '  -3.46%  '
Set-TextValue $ws.Range("B31") 'PancakeSwap'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D31") '1.551'
Set-TextValue $ws.Range("E31") '  -2.59%  '
Set-TextValue $ws.Range("B32") 'Filecoin'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D32") '4.554'
Set-TextValue $ws.Range("E32") '  -4.75%  '
Set-TextValue $ws.Range("D33") '4.365'
Set-TextValue $ws.Range("E33") '  -3.52%  '
Set-TextValue $ws.Range("D34") '0.04846'
Set-TextValue $ws.Range("E34") '  -3.14%  '
Set-TextValue $ws.Range("D35") '0.7511'
Set-TextValue $ws.Range("E35") '  -1.92%  '
Set-TextValue $ws.Range("D36") '1.158'
Set-TextValue $ws.Range("E36") '  -0.65%  '
Set-TextValue $ws.Range("D37") '0.9995'
Set-TextValue $ws.Range("E37") '  -0.09%  '
Set-TextValue $ws.Range("D38") '2.724'
Set-TextValue $ws.Range("E38") '  -0.33%  '
Set-TextValue $ws.Range("D39") '0.01991'
Set-TextValue $ws.Range("E39") '  -2.97%  '
Set-TextValue $ws.Range("D40") '2.660'
Set-TextValue $ws.Range("E40") '  -1.96%  '
Set-TextValue $ws.Range("D41") '78.54'
Set-TextValue $ws.Range("E41") '  +7.43%  '
Set-TextValue $ws.Range("D42") '6.477'
Set-TextValue $ws.Range("E42") '  +0.94%  '
Set-TextValue $ws.Range("D43") '2.095'
Set-TextValue $ws.Range("E43") '  -1.63%  '
Set-TextValue $ws.Range("D44") '0.9065'
Set-TextValue $ws.Range("E44") '  +2.42%  '
Set-TextValue $ws.Range("D45") '108.13'
Set-TextValue $ws.Range("E45") '  -1.17%  '
Set-TextValue $ws.Range("E46") '  -1.42%  '
Set-TextValue $ws.Range("D47") '0.9986'
Set-TextValue $ws.Range("E47") '  -0.19%  '
Set-TextValue $ws.Range("D48") '7.752'
Set-TextValue $ws.Range("E48") '  +3.82%  '
Set-TextValue $ws.Range("D49") '988.01'
Set-TextValue $ws.Range("E49") '  -3.17%  '
Set-TextValue $ws.Range("E50") '  -1.88%  '
Set-TextValue $ws.Range("D51") '9.287'
Set-TextValue $ws.Range("E51") '  -0.90%  '
